# Updated cryptos list with latest price/volume data.
# Row 36-38 are also reordered: Monero, ImmutableX, NEARProtocol.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.895.54"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "2.418.66"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'569.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").Value = "'139.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").Value = "2.403.28"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").Value = "'0.107"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").Value = "'0.159"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "'5.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").Value = "'26.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "2.830.18"
$ws.Range("D17").Value = "60.835.63"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "2.419.75"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "'7.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.61%  "
$ws.Range("D20").Value = "'10.63"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "'322.70"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'4.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").Value = "'6.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'1.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.73%  "
$ws.Range("D26").Value = "'64.77"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").Value = "'580.46"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").Value = "'8.28"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -10.03%  "
$ws.Range("D29").Value = "2.540.81"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").Value = "'7.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "'1.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.22%  "
$ws.Range("D33").Value = "'1.84"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'152.44"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.40"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.60"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("D39").Value = "'0.367"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").Value = "'18.25"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").Value = "'5.13"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'1.67"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'41.19"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("D45").Value = "'2.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.82%  "
$ws.Range("D46").Value = "0.0₆0281"
$ws.Range("E46").Value = "  +5.13%  "
$ws.Range("D47").Value = "'142.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'3.51"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("D49").Value = "'0.587"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("D50").Value = "'19.36"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'0.0502"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.61%  "
